# Fix: Elective lecture and tutorial scheduling
# Updates the lecture/tutorial room codes (columns D and E) for rows 20-35
# on the Regular_Timetable, PreMid_Timetable and PostMid_Timetable sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Regular_Timetable", "PreMid_Timetable", "PostMid_Timetable")

$rowData = @{
    20 = @{ D = "Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]"; E = "Tue 14:30-15:30 [C101]" }
    21 = @{ D = "Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]"; E = "Tue 14:30-15:30 [C102]" }
    22 = @{ D = "Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]"; E = "Tue 14:30-15:30 [C104]" }
    23 = @{ D = "Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]"; E = "Wed 14:30-15:30 [C101]" }
    24 = @{ D = "Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]"; E = "Wed 14:30-15:30 [C102]" }
    25 = @{ D = "Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]"; E = "Wed 14:30-15:30 [C104]" }
    26 = @{ D = "Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]"; E = "Wed 14:30-15:30 [C202]" }
    27 = @{ D = "Mon 10:30-12:00 [C101], Wed 10:30-12:00 [C101]"; E = "Thu 14:30-15:30 [C101]" }
    28 = @{ D = "Mon 10:30-12:00 [C102], Wed 10:30-12:00 [C102]"; E = "Thu 14:30-15:30 [C102]" }
    29 = @{ D = "Mon 10:30-12:00 [C104], Wed 10:30-12:00 [C104]"; E = "Thu 14:30-15:30 [C104]" }
    30 = @{ D = "Mon 10:30-12:00 [C202], Wed 10:30-12:00 [C202]"; E = "Thu 14:30-15:30 [C202]" }
    31 = @{ D = "Mon 10:30-12:00 [C203], Wed 10:30-12:00 [C203]"; E = "Thu 14:30-15:30 [C203]" }
    32 = @{ D = "Tue 15:30-17:00 [C101], Thu 15:30-17:00 [C101]"; E = "Fri 14:30-15:30 [C101]" }
    33 = @{ D = "Tue 15:30-17:00 [C102], Thu 15:30-17:00 [C102]"; E = "Fri 14:30-15:30 [C102]" }
    34 = @{ D = "Tue 15:30-17:00 [C104], Thu 15:30-17:00 [C104]"; E = "Fri 14:30-15:30 [C104]" }
    35 = @{ D = "Tue 15:30-17:00 [C202], Thu 15:30-17:00 [C202]"; E = "Fri 14:30-15:30 [C202]" }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $rowData.Keys) {
        $vals = $rowData[$row]
        $ws.Range("D$row").Value = $vals.D
        $ws.Range("E$row").Value = $vals.E
    }
}
